$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename tc01 -> testCase01 (A2 currently holds "tc01")
$ws.Range("A2").Value = "testCase01"

# Add new rows 3-8 for testCase02..testCase07
$testCases = @("testCase02", "testCase03", "testCase04", "testCase05", "testCase06", "testCase07")

$row = 3
foreach ($tc in $testCases) {
    $ws.Cells.Item($row, 1).Value = $tc
    $ws.Cells.Item($row, 2).Value = "openBrowser"
    $ws.Cells.Item($row, 3).Value = "launchUrl"
    $ws.Cells.Item($row, 4).Value = "login"
    $ws.Cells.Item($row, 6).Value = "closeBrowser"
    $row++
}

# Update the active selection to D10 as in the diff
$ws.Range("D10").Select() | Out-Null
